# Apply the permuted dataset values for rows 2-10 (columns A,B,D,E,F,G,H,I,Q,R,AO).
# Each row keeps its static/context columns (C,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AI,AT,AW,AX,AY)
# and only the species-record fields move between rows, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 111609174
$ws.Cells.Item(2,2).Value = 90854
$ws.Cells.Item(2,4).Value = "NT"
$ws.Cells.Item(2,5).Value = 2079
$ws.Cells.Item(2,6).Value = "Nordtagging"
$ws.Cells.Item(2,7).Value = "Odonticium romellii"
$ws.Cells.Item(2,8).Value = "(S.Lundell) Parmasto"
$ws.Cells.Item(2,9).Value = ""
$ws.Cells.Item(2,17).Value = 514788.8674634451
$ws.Cells.Item(2,18).Value = 6925250.666874606
$ws.Cells.Item(2,41).Value = "mossig silverlåga av tall"

# Row 3
$ws.Cells.Item(3,1).Value = 111609168
$ws.Cells.Item(3,2).Value = 77597
$ws.Cells.Item(3,4).Value = "NT"
$ws.Cells.Item(3,5).Value = 864
$ws.Cells.Item(3,6).Value = "Knottrig blåslav"
$ws.Cells.Item(3,7).Value = "Hypogymnia bitteri"
$ws.Cells.Item(3,8).Value = "(Lynge) Ahti"
$ws.Cells.Item(3,9).Value = ""
$ws.Cells.Item(3,17).Value = 515085.0087401169
$ws.Cells.Item(3,18).Value = 6925147.4056778
$ws.Cells.Item(3,41).Value = "tall"

# Row 4
$ws.Cells.Item(4,1).Value = 111609172
$ws.Cells.Item(4,2).Value = 77268
$ws.Cells.Item(4,4).Value = "NT"
$ws.Cells.Item(4,5).Value = 228912
$ws.Cells.Item(4,6).Value = "Mörk kolflarnlav"
$ws.Cells.Item(4,7).Value = "Carbonicola myrmecina"
$ws.Cells.Item(4,8).Value = "(Ach.) Bendiksby & Timdal"
$ws.Cells.Item(4,9).Value = ""
$ws.Cells.Item(4,17).Value = 514955.9350709137
$ws.Cells.Item(4,18).Value = 6925302.779521272
$ws.Cells.Item(4,41).Value = "brandstubbe"

# Row 5
$ws.Cells.Item(5,1).Value = 111609176
$ws.Cells.Item(5,2).Value = 96348
$ws.Cells.Item(5,4).Value = "VU"
$ws.Cells.Item(5,5).Value = 220787
$ws.Cells.Item(5,6).Value = "Knärot"
$ws.Cells.Item(5,7).Value = "Goodyera repens"
$ws.Cells.Item(5,8).Value = "(L.) R. Br."
$ws.Cells.Item(5,9).NumberFormat = "@"
$ws.Cells.Item(5,9).Value = "30"
$ws.Cells.Item(5,17).Value = 514875.4249116365
$ws.Cells.Item(5,18).Value = 6924962.732657854
$ws.Cells.Item(5,41).Value = ""

# Row 6
$ws.Cells.Item(6,1).Value = 111609170
$ws.Cells.Item(6,2).Value = 96348
$ws.Cells.Item(6,4).Value = "VU"
$ws.Cells.Item(6,5).Value = 220787
$ws.Cells.Item(6,6).Value = "Knärot"
$ws.Cells.Item(6,7).Value = "Goodyera repens"
$ws.Cells.Item(6,8).Value = "(L.) R. Br."
$ws.Cells.Item(6,9).NumberFormat = "@"
$ws.Cells.Item(6,9).Value = "3"
$ws.Cells.Item(6,17).Value = 515035.9338400747
$ws.Cells.Item(6,18).Value = 6925238.814452391
$ws.Cells.Item(6,41).Value = ""

# Row 7
$ws.Cells.Item(7,1).Value = 111609175
$ws.Cells.Item(7,2).Value = 77268
$ws.Cells.Item(7,4).Value = "NT"
$ws.Cells.Item(7,5).Value = 228912
$ws.Cells.Item(7,6).Value = "Mörk kolflarnlav"
$ws.Cells.Item(7,7).Value = "Carbonicola myrmecina"
$ws.Cells.Item(7,8).Value = "(Ach.) Bendiksby & Timdal"
$ws.Cells.Item(7,9).Value = ""
$ws.Cells.Item(7,17).Value = 514769.8196280882
$ws.Cells.Item(7,18).Value = 6925156.6384242
$ws.Cells.Item(7,41).Value = "brandstubbe"

# Row 8
$ws.Cells.Item(8,1).Value = 111609167
$ws.Cells.Item(8,2).Value = 77186
$ws.Cells.Item(8,4).Value = "NT"
$ws.Cells.Item(8,5).Value = 353
$ws.Cells.Item(8,6).Value = "Dvärgbägarlav"
$ws.Cells.Item(8,7).Value = "Cladonia parasitica"
$ws.Cells.Item(8,8).Value = "(Hoffm.) Hoffm."
$ws.Cells.Item(8,9).Value = ""
$ws.Cells.Item(8,17).Value = 515051.1877758073
$ws.Cells.Item(8,18).Value = 6925144.938876954
$ws.Cells.Item(8,41).Value = "silverlåga av tall"

# Row 9
$ws.Cells.Item(9,1).Value = 111609169
$ws.Cells.Item(9,2).Value = 96348
$ws.Cells.Item(9,4).Value = "VU"
$ws.Cells.Item(9,5).Value = 220787
$ws.Cells.Item(9,6).Value = "Knärot"
$ws.Cells.Item(9,7).Value = "Goodyera repens"
$ws.Cells.Item(9,8).Value = "(L.) R. Br."
$ws.Cells.Item(9,9).NumberFormat = "@"
$ws.Cells.Item(9,9).Value = "4"
$ws.Cells.Item(9,17).Value = 515078.8479096842
$ws.Cells.Item(9,18).Value = 6925177.45879681
$ws.Cells.Item(9,41).Value = ""

# Row 10
$ws.Cells.Item(10,1).Value = 111609173
$ws.Cells.Item(10,2).Value = 96348
$ws.Cells.Item(10,4).Value = "VU"
$ws.Cells.Item(10,5).Value = 220787
$ws.Cells.Item(10,6).Value = "Knärot"
$ws.Cells.Item(10,7).Value = "Goodyera repens"
$ws.Cells.Item(10,8).Value = "(L.) R. Br."
$ws.Cells.Item(10,9).NumberFormat = "@"
$ws.Cells.Item(10,9).Value = "7"
$ws.Cells.Item(10,17).Value = 514934.1293421969
$ws.Cells.Item(10,18).Value = 6925308.234934391
$ws.Cells.Item(10,41).Value = ""
